$d = $word.ActiveDocument

$xml = $d.WordOpenXML

# The BTec_Logo-Orange pictures (in the headers) currently carry the
# displayed name "image2.jpg" and need to become "image1.jpg".
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

# The Pearson logo pictures (in the footers) currently carry the
# displayed name "image1.png" and need to become "image2.png".
$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')

$d.WordOpenXML = $xml
